$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    # Force the value to be stored as text even when it looks numeric
    # (organismID catalog numbers), then drop back to the default/unstyled
    # cell format so no stray style index is left attached to the cell.
    $r = $ws.Range($rangeAddr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Insert a new row before row 32, pushing the current rows 32 and 33 down to 33 and 34.
$ws.Rows("32:32").Insert()

# New row 32 — record for "Tecoma peroba Record" (same underlying collection event as
# row 31, but its own coordinates/issues/organismID).
$ws.Range("A32").Value = "Tecoma peroba Record"
$ws.Range("B32").Value = "Paratecoma peroba (Record) Kuhlm."
$ws.Range("C32").Value = -10.8339
$ws.Range("D32").Value = -52.8731
$ws.Range("E32").Value = "gass84,ambcol,inmafu"
$ws.Range("F32").Value = "PRESERVED_SPECIMEN"
$ws.Range("G32").Value = "PRESENT"
$ws.Range("H32").Value = "The Field Museum of Natural History"
$ws.Range("I32").Value = "vTypes"
$ws.Range("J32").Value = "H. N. Whitford & F. Silveira"
Set-TextValue "L32" "271506"

# New column L: organismID, with the same (bold/centered) formatting as the other header cells.
$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial(-4122)
$ws.Range("L1").Value = "organismID"

# organismID values for the existing rows 29-31.
Set-TextValue "L29" "3736349"
Set-TextValue "L30" "323929"
Set-TextValue "L31" "271507"

$excel.CutCopyMode = 0
